# Generate Report for Handback
# Updates the "zh-cn" and "de-de" detail sheets with the handback info for
# the ea91ca86-d878-4343-a738-8868285667ae entry (row 8), plus widens the
# "Error Detail" column so the new message is readable.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/88c2e0e3e5dd553451b671bfbfe532c8fbc823dd/e2e/ea91ca86-d878-4343-a738-8868285667ae.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e495ad49933fdda1538f584b2c9d2ba3b303ef6b/e2e/ea91ca86-d878-4343-a738-8868285667ae.md."
$latestHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e495ad49933fdda1538f584b2c9d2ba3b303ef6b/e2e/ea91ca86-d878-4343-a738-8868285667ae.md"
$handbackDisplay = "ea91ca86-d878-4343-a738-8868285667ae.md"

function Update-LocaleSheet($sheetName, $xlfFileName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Latest Target File (I8): link back to the handback markdown file, same
    # as column A, now that a handback has actually arrived for this row.
    $ws.Hyperlinks.Add($ws.Range("I8"), $latestHandbackUrl, "", "", $handbackDisplay)
    $ws.Range("I8").ClearFormats()
    $ws.Range("I8").Font.Underline = $true
    $ws.Range("I8").Font.Color = 15570276

    # Latest Handback File (J8)
    $ws.Range("J8").Value = $xlfFileName

    # Latest Handback DateTime (K8)
    $ws.Range("K8").Value = $handbackDateTime

    # Error Detail (P8)
    $ws.Range("P8").Value = $errorMessage

    # Widen the Error Detail column (P) so the long message is legible.
    $ws.Range("P1").ColumnWidth = 39.1666666666667
}

Update-LocaleSheet "zh-cn" "ea91ca86-d878-4343-a738-8868285667ae.1ac22fbc857a6385ba5bede3ef45ecf09d05062c.zh-cn.xlf" "2016-09-03 10:46:32"
Update-LocaleSheet "de-de" "ea91ca86-d878-4343-a738-8868285667ae.1ac22fbc857a6385ba5bede3ef45ecf09d05062c.de-de.xlf" "2016-09-03 10:46:39"
